$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure columns D, E, G keep their values as literal text (not auto-converted
# to numbers/percentages/dates by Excel) by marking them as Text format first.
$ws.Range("D2:E51").NumberFormat = "@"
$ws.Range("G2:G51").NumberFormat = "@"

$ws.Range("D2").Value = '246.03'
$ws.Range("E2").Value = '-0.17%'
$ws.Range("G2").Value = '23'
$ws.Range("D3").Value = '28.25'
$ws.Range("E3").Value = '-3.65%'
$ws.Range("G3").Value = '23'
$ws.Range("D4").Value = '5.293'
$ws.Range("E4").Value = '1.85%'
$ws.Range("G4").Value = '23'
$ws.Range("D5").Value = '0.05712'
$ws.Range("E5").Value = '-0.39%'
$ws.Range("G5").Value = '23'
$ws.Range("D6").Value = '6.651'
$ws.Range("E6").Value = '1.32%'
$ws.Range("G6").Value = '23'
$ws.Range("D7").Value = '3.215'
$ws.Range("E7").Value = '3.53%'
$ws.Range("G7").Value = '23'
$ws.Range("D8").Value = '0.8636'
$ws.Range("E8").Value = '0.57%'
$ws.Range("G8").Value = '23'
$ws.Range("D9").Value = '0.8906'
$ws.Range("E9").Value = '3.10%'
$ws.Range("G9").Value = '23'
$ws.Range("D10").Value = '0.1387'
$ws.Range("E10").Value = '1.75%'
$ws.Range("G10").Value = '23'
$ws.Range("D11").Value = '0.07100'
$ws.Range("E11").Value = '0.35%'
$ws.Range("G11").Value = '23'
$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D12").Value = '0.03125'
$ws.Range("E12").Value = '2.47%'
$ws.Range("G12").Value = '23'
$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D13").Value = '0.09234'
$ws.Range("E13").Value = '-1.52%'
$ws.Range("G13").Value = '23'
$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D14").Value = '0.001527'
$ws.Range("E14").Value = '-0.54%'
$ws.Range("G14").Value = '23'
$ws.Range("B15").Value = 'One'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D15").Value = '0.0005987'
$ws.Range("E15").Value = '-0.76%'
$ws.Range("G15").Value = '23'
$ws.Range("B16").Value = 'TigerCash'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D16").Value = '0.006036'
$ws.Range("E16").Value = '-0.56%'
$ws.Range("G16").Value = '23'
$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D17").Value = '3.495'
$ws.Range("E17").Value = '0.09%'
$ws.Range("G17").Value = '23'
$ws.Range("B18").Value = 'BTSEToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D18").Value = '2.172'
$ws.Range("E18").Value = '-4.77%'
$ws.Range("G18").Value = '23'
$ws.Range("B19").Value = 'BitpandaEcosystemToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D19").Value = '0.3167'
$ws.Range("E19").Value = '-1.05%'
$ws.Range("G19").Value = '23'
$ws.Range("B20").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C20").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D20").Value = '0.03333'
$ws.Range("E20").Value = '0.49%'
$ws.Range("G20").Value = '23'
$ws.Range("E21").Value = '1.47%'
$ws.Range("G21").Value = '23'
$ws.Range("D22").Value = '3.492'
$ws.Range("E22").Value = '0.73%'
$ws.Range("G22").Value = '23'
$ws.Range("D23").Value = '0.04086'
$ws.Range("E23").Value = '-1.62%'
$ws.Range("G23").Value = '23'
$ws.Range("G24").Value = '23'
$ws.Range("D25").Value = '0.001218'
$ws.Range("E25").Value = '-0.58%'
$ws.Range("G25").Value = '23'
$ws.Range("D26").Value = '0.004166'
$ws.Range("E26").Value = '-16.61%'
$ws.Range("G26").Value = '23'
$ws.Range("D27").Value = '0.0001199'
$ws.Range("G27").Value = '23'
$ws.Range("D28").Value = '0.0001444'
$ws.Range("G28").Value = '23'
$ws.Range("G29").Value = '23'
$ws.Range("G30").Value = '23'
$ws.Range("G31").Value = '23'
$ws.Range("G32").Value = '23'
$ws.Range("G33").Value = '23'
$ws.Range("G34").Value = '23'
$ws.Range("G35").Value = '23'
$ws.Range("G36").Value = '23'
$ws.Range("G37").Value = '23'
$ws.Range("G38").Value = '23'
$ws.Range("G39").Value = '23'
$ws.Range("D40").Value = '0.03791'
$ws.Range("E40").Value = '1.05%'
$ws.Range("G40").Value = '23'
$ws.Range("E41").Value = '-0.33%'
$ws.Range("G41").Value = '23'
$ws.Range("D42").Value = '0.002407'
$ws.Range("E42").Value = '14.62%'
$ws.Range("G42").Value = '23'
$ws.Range("D43").Value = '0.002947'
$ws.Range("E43").Value = '-48.94%'
$ws.Range("G43").Value = '23'
$ws.Range("D44").Value = '0.009456'
$ws.Range("E44").Value = '11.78%'
$ws.Range("G44").Value = '23'
$ws.Range("D45").Value = '0.00005269'
$ws.Range("E45").Value = '-0.20%'
$ws.Range("G45").Value = '23'
$ws.Range("E46").Value = '-0.06%'
$ws.Range("G46").Value = '23'
$ws.Range("D47").Value = '0.08905'
$ws.Range("E47").Value = '56.17%'
$ws.Range("G47").Value = '23'
$ws.Range("D48").Value = '0.002255'
$ws.Range("E48").Value = '-0.86%'
$ws.Range("G48").Value = '23'
$ws.Range("E49").Value = '-0.06%'
$ws.Range("G49").Value = '23'
$ws.Range("E50").Value = '-0.06%'
$ws.Range("G50").Value = '23'
$ws.Range("G51").Value = '23'
